# Apply the "added lidar module" OutListParameters.xlsx edit:
#   - InflowWind sheet: new "Wind Sensor Measurements" category (row 30)
#     followed by 5 new output channels WindMeas1..WindMeas5 (rows 31-35).
#   - Instructions sheet: COUNTA formula recalculates (27 -> 32) and the
#     active selection moves from D9 to D8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InflowWind")
$instr = $wb.Worksheets.Item("Instructions")

# --- Row 30: new category header "Wind Sensor Measurements" ---------------
# Copy the style of the previous category row (row 2: A="s3", C/D/E="s4")
# onto row 30, then overwrite the text.
$ws.Range("A2").Copy($ws.Range("A30"))
$ws.Range("C2:E2").Copy($ws.Range("C30:E30"))
$ws.Range("A30").Value = "Wind Sensor Measurements"
# The old placeholder cell on F30 (s=5, no content) is not part of the new
# category row - remove it entirely (contents + formatting).
$ws.Range("F30").Clear()

# --- Rows 31-35: WindMeas1..WindMeas5 output channels ----------------------
# Seed style from the last existing data row (29), which already carries the
# B/D/E/F formatting (s5/s6/s6/s5) used throughout the table.
$ws.Range("B29:F29").Copy($ws.Range("B31:F31"))
$ws.Range("B29:F29").Copy($ws.Range("B32:F32"))
$ws.Range("B29:F29").Copy($ws.Range("B33:F33"))
$ws.Range("B29:F29").Copy($ws.Range("B34:F34"))
$ws.Range("B29:F29").Copy($ws.Range("B35:F35"))

$names = @("WindMeas1", "WindMeas2", "WindMeas3", "WindMeas4", "WindMeas5")
$descs = @(
    "Wind measurement at sensor 1",
    "Wind measurement at sensor 2",
    "Wind measurement at sensor 3",
    "Wind measurement at sensor 4",
    "Wind measurement at sensor 5"
)
$criteria = @(
    "p%lidar%SensorType == SensorType_None",
    "p%lidar%NumPulseGate < 2",
    "p%lidar%NumPulseGate < 3",
    "p%lidar%NumPulseGate < 4",
    "p%lidar%NumPulseGate < 5"
)

for ($i = 0; $i -lt 5; $i++) {
    $r = 31 + $i
    $ws.Cells.Item($r, 2).Value = $names[$i]
    $ws.Cells.Item($r, 3).Value = ""
    $ws.Cells.Item($r, 4).Value = $descs[$i]
    $ws.Cells.Item($r, 5).Value = "Defined by sensor"
    $ws.Cells.Item($r, 6).Value = "(m/s)"
    $ws.Cells.Item($r, 7).Value = $criteria[$i]
}

# C column on the new data rows should stay empty (clear the copied style).
$ws.Range("C31:C35").Value = ""

# --- Sheet selections, matching where the author ended up -----------------
# Instructions' own cursor moved from D9 to D8 ...
$instr.Range("D8").Select()
# ... but InflowWind is the sheet left active/visible (tabSelected) with its
# cursor on the first newly-added cell.
$ws.Select()
$ws.Range("B31").Select()

$wb.Save()
